$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1168")

# Corrected "Section" (column H) values that were previously input incorrectly.
$updates = @{
    222 = 2; 223 = 2
    234 = 2; 235 = 2; 236 = 2
    238 = 3; 239 = 3; 240 = 2; 241 = 3; 242 = 2; 243 = 4; 244 = 3; 245 = 4
    246 = 3; 247 = 3; 248 = 3; 249 = 4; 250 = 3; 251 = 2; 252 = 2; 253 = 2
    255 = 2; 256 = 3; 257 = 3; 258 = 3; 259 = 2; 260 = 3; 261 = 2; 262 = 2
    263 = 2; 264 = 2; 266 = 3; 267 = 3; 268 = 3; 269 = 3; 270 = 2; 271 = 4
    272 = 3; 273 = 4; 274 = 2; 275 = 3; 276 = 4; 277 = 3
    284 = 2
    286 = 2; 287 = 2; 288 = 2
    290 = 2
    295 = 2; 296 = 2; 297 = 2; 298 = 2
    303 = 2; 304 = 2; 305 = 2
    307 = 2
    312 = 2
    315 = 2; 316 = 2; 317 = 2; 318 = 2; 319 = 4; 320 = 4; 321 = 3; 322 = 3
    324 = 3
    326 = 3; 327 = 3
    329 = 2
    332 = 2
    338 = 2; 339 = 2; 340 = 2; 341 = 2
    355 = 2; 356 = 2; 357 = 2
    364 = 2; 365 = 2
    367 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("H$row").Value = $updates[$row]
}

# Restore the view to the top of the sheet and move the active selection,
# matching the state left behind after the corrections were made.
$ws.Activate()
$ws.Range("H224").Select()
$excel.ActiveWindow.ScrollRow = 1
